$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin prices (Price column, D) for the Dec 25 2022 symbol-list refresh.
# Source cells store prices as literal text (t="inlineStr"), not numbers, so we
# force a Text number format while writing the new value and then drop the format
# back to the sheet's default style (the value itself stays text either way since
# it was written while the cell was Text-formatted).
$priceUpdates = [ordered]@{
    "D2" = "245.32"
    "D3" = "23.12"
    "D4" = "5.401"
    "D5" = "0.06047"
    "D6" = "3.400"
    "D7" = "0.8064"
    "D8" = "0.9355"
    "D9" = "0.1429"
    "D10" = "0.07441"
    "D11" = "0.03331"
    "D12" = "0.03059"
    "D13" = "4.008"
    "D14" = "0.09366"
    "D16" = "0.04822"
    "D17" = "0.0005944"
    "D18" = "0.005489"
    "D20" = "0.0009865"
    "D21" = "0.00008705"
    "D22" = "3.658"
    "D23" = "6.434"
    "D27" = "0.0002902"
    "D40" = "0.03984"
    "D41" = "0.006373"
    "D43" = "0.002711"
    "D44" = "0.006296"
    "D45" = "0.00005253"
    "D46" = "0.00000000751"
    "D47" = "0.0005804"
    "D48" = "1.091"
    "D49" = "0.002236"
    "D51" = "0.01011"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}
